$d = $word.ActiveDocument

$p2 = $d.Paragraphs.Item(2)
$rng2 = $p2.Range
$xml2 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="28B13BAB" w14:textId="77777777" w:rsidR="00655392" w:rsidRDefault="003D6324"><w:pPr><w:widowControl w:val="0"/><w:tabs><w:tab w:val="right" w:pos="12150"/></w:tabs><w:jc w:val="both"/><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:sz w:val="40"/><w:szCs w:val="40"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">Project Name:  Project 1:  </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t>CompuVote</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:b/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr><w:t xml:space="preserve">                                                                                                    Team #19</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng2.InsertXML($xml2)
Write-Output "Done paragraph 2"

$p9 = $d.Paragraphs.Item(9)
$rng9 = $p9.Range
$xml9 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="5C583570" w14:textId="4F2C78EB" w:rsidR="00655392" w:rsidRDefault="003D6324"><w:pPr><w:widowControl w:val="0"/></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t>Name(s) of Testers:</w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="006670B7"><w:t xml:space="preserve">Aaron </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="006670B7"><w:t>Kandikatla</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng9.InsertXML($xml9)
Write-Output "Done paragraph 9"

$p17 = $d.Paragraphs.Item(17)
$rng17 = $p17.Range
$xml17 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="357578E5" w14:textId="411CD4AC" w:rsidR="002E5151" w:rsidRPr="002E5151" w:rsidRDefault="002E5151" w:rsidP="002E5151"><w:pPr><w:pStyle w:val="NormalWeb"/></w:pPr><w:r><w:t>Test file: Project1/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>src</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/test/org/team19/OpenPartyListSystemTest.java</w:t></w:r><w:r><w:br/><w:t xml:space="preserve">Test method: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>t</w:t></w:r><w:r w:rsidR="009E6689"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>est</w:t></w:r><w:r w:rsidR="004F1D24"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>PrintSummaryTable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:br/><w:t xml:space="preserve">Method/constructor being tested: </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidR="004F1D24"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/></w:rPr><w:t>printSummaryTable</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> from Project1/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>src</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/main/org/team19/OpenPartyListSystem.java</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng17.InsertXML($xml17)
Write-Output "Done paragraph 17"

$p28 = $d.Paragraphs.Item(28)
$rng28 = $p28.Range
$xml28 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="0E0C7E14" w14:textId="12EB3CE4" w:rsidR="00655392" w:rsidRDefault="004F122F"><w:pPr><w:widowControl w:val="0"/></w:pPr><w:bookmarkStart w:id="2" w:name="OLE_LINK33"/><w:bookmarkStart w:id="3" w:name="OLE_LINK34"/><w:r><w:t xml:space="preserve">There exists a file called </w:t></w:r><w:bookmarkStart w:id="4" w:name="OLE_LINK24"/><w:bookmarkStart w:id="5" w:name="OLE_LINK25"/><w:bookmarkStart w:id="6" w:name="OLE_LINK63"/><w:r w:rsidRPr="004F122F"><w:t>test_print_summary_table_expected</w:t></w:r><w:r w:rsidRPr="00FB219E"><w:t>.txt</w:t></w:r><w:bookmarkEnd w:id="4"/><w:bookmarkEnd w:id="5"/><w:bookmarkEnd w:id="6"/><w:r w:rsidRPr="00FB219E"><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>in testing/test-</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>recources</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t>/</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>openPartyListSystemTest</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve">/ which represents the expected </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>output</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> and the system is able to open and read the file.</w:t></w:r><w:bookmarkEnd w:id="2"/><w:bookmarkEnd w:id="3"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng28.InsertXML($xml28)
Write-Output "Done paragraph 28"

$p75 = $d.Paragraphs.Item(75)
$rng75 = $p75.Range
$xml75 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="217F2196" w14:textId="1072B651" w:rsidR="00655392" w:rsidRDefault="004F122F" w:rsidP="00A976C2"><w:pPr><w:widowControl w:val="0"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:bookmarkStart w:id="9" w:name="OLE_LINK11"/><w:bookmarkStart w:id="10" w:name="OLE_LINK12"/><w:bookmarkStart w:id="11" w:name="OLE_LINK32"/><w:r w:rsidRPr="004F122F"><w:lastRenderedPageBreak/><w:t>test_print_summary_table_expected</w:t></w:r><w:r w:rsidRPr="00FB219E"><w:t>.txt</w:t></w:r><w:r><w:t xml:space="preserve"> matches the file written by </w:t></w:r><w:bookmarkEnd w:id="9"/><w:bookmarkEnd w:id="10"/><w:bookmarkEnd w:id="11"/><w:proofErr w:type="spellStart"/><w:r><w:t>p</w:t></w:r><w:r w:rsidRPr="004F122F"><w:t>rintSummaryTable</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng75.InsertXML($xml75)
Write-Output "Done paragraph 75"

$p76 = $d.Paragraphs.Item(76)
$rng76 = $p76.Range
$xml76 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="25384A84" w14:textId="5BC1ADBA" w:rsidR="00655392" w:rsidRDefault="004F122F" w:rsidP="00A976C2"><w:pPr><w:widowControl w:val="0"/><w:rPr><w:sz w:val="22"/><w:szCs w:val="22"/></w:rPr></w:pPr><w:r w:rsidRPr="004F122F"><w:t>test_print_summary_table_expected</w:t></w:r><w:r w:rsidRPr="00FB219E"><w:t>.txt</w:t></w:r><w:r><w:t xml:space="preserve"> matches the file written by </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>p</w:t></w:r><w:r w:rsidRPr="004F122F"><w:t>rintSummaryTable</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng76.InsertXML($xml76)
Write-Output "Done paragraph 76"

$p81 = $d.Paragraphs.Item(81)
$rng81 = $p81.Range
$xml81 = '<?xml version="1.0"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships"><w:body><w:p w14:paraId="53C95BE5" w14:textId="4D9B8469" w:rsidR="00655392" w:rsidRPr="00D6299E" w:rsidRDefault="003D6324"><w:pPr><w:pBdr><w:top w:val="nil"/><w:left w:val="nil"/><w:bottom w:val="nil"/><w:right w:val="nil"/><w:between w:val="nil"/></w:pBdr><w:spacing w:after="140" w:line="288" w:lineRule="auto"/></w:pPr><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wpg"><w:drawing><wp:anchor distT="0" distB="0" distL="0" distR="0" simplePos="0" relativeHeight="251658240" behindDoc="0" locked="0" layoutInCell="1" hidden="0" allowOverlap="1" wp14:anchorId="6ED74B1D" wp14:editId="54697917"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>0</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>-165099</wp:posOffset></wp:positionV><wp:extent cx="8363585" cy="12700"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapSquare wrapText="bothSides" distT="0" distB="0" distL="0" distR="0"/><wp:docPr id="1" name="Straight Arrow Connector 1"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvCnPr/><wps:spPr><a:xfrm><a:off x="1164600" y="3780000"/><a:ext cx="8362800" cy="0"/></a:xfrm><a:prstGeom prst="straightConnector1"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525" cap="flat" cmpd="sng"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:prstDash val="solid"/><a:miter lim="8000"/><a:headEnd type="none" w="sm" len="sm"/><a:tailEnd type="none" w="sm" len="sm"/></a:ln></wps:spPr><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback xmlns:sl="http://schemas.openxmlformats.org/schemaLibrary/2006/main" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:c="http://schemas.openxmlformats.org/drawingml/2006/chart" xmlns:lc="http://schemas.openxmlformats.org/drawingml/2006/lockedCanvas" xmlns:dgm="http://schemas.openxmlformats.org/drawingml/2006/diagram"><w:drawing><wp:anchor allowOverlap="1" behindDoc="0" distB="0" distT="0" distL="0" distR="0" hidden="0" layoutInCell="1" locked="0" relativeHeight="0" simplePos="0"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>0</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>-165099</wp:posOffset></wp:positionV><wp:extent cx="8363585" cy="12700"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:wrapSquare wrapText="bothSides" distB="0" distT="0" distL="0" distR="0"/><wp:docPr id="1" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="image1.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId7"/><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="8363585" cy="12700"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Fallback></mc:AlternateContent></w:r><w:r><w:rPr><w:noProof/></w:rPr><mc:AlternateContent><mc:Choice Requires="wpg"><w:drawing><wp:anchor distT="0" distB="0" distL="0" distR="0" simplePos="0" relativeHeight="251659264" behindDoc="0" locked="0" layoutInCell="1" hidden="0" allowOverlap="1" wp14:anchorId="6D575460" wp14:editId="4A4CCBBA"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>0</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>533400</wp:posOffset></wp:positionV><wp:extent cx="8363585" cy="12700"/><wp:effectExtent l="0" t="0" r="0" b="0"/><wp:wrapSquare wrapText="bothSides" distT="0" distB="0" distL="0" distR="0"/><wp:docPr id="2" name="Straight Arrow Connector 2"/><wp:cNvGraphicFramePr/><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.microsoft.com/office/word/2010/wordprocessingShape"><wps:wsp><wps:cNvCnPr/><wps:spPr><a:xfrm><a:off x="1164600" y="3780000"/><a:ext cx="8362800" cy="0"/></a:xfrm><a:prstGeom prst="straightConnector1"><a:avLst/></a:prstGeom><a:noFill/><a:ln w="9525" cap="flat" cmpd="sng"><a:solidFill><a:srgbClr val="000000"/></a:solidFill><a:prstDash val="solid"/><a:miter lim="8000"/><a:headEnd type="none" w="sm" len="sm"/><a:tailEnd type="none" w="sm" len="sm"/></a:ln></wps:spPr><wps:bodyPr/></wps:wsp></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Choice><mc:Fallback xmlns:sl="http://schemas.openxmlformats.org/schemaLibrary/2006/main" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:c="http://schemas.openxmlformats.org/drawingml/2006/chart" xmlns:lc="http://schemas.openxmlformats.org/drawingml/2006/lockedCanvas" xmlns:dgm="http://schemas.openxmlformats.org/drawingml/2006/diagram"><w:drawing><wp:anchor allowOverlap="1" behindDoc="0" distB="0" distT="0" distL="0" distR="0" hidden="0" layoutInCell="1" locked="0" relativeHeight="0" simplePos="0"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="column"><wp:posOffset>0</wp:posOffset></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>533400</wp:posOffset></wp:positionV><wp:extent cx="8363585" cy="12700"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:wrapSquare wrapText="bothSides" distB="0" distT="0" distL="0" distR="0"/><wp:docPr id="2" name="image2.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr id="0" name="image2.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId8"/><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="8363585" cy="12700"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:anchor></w:drawing></mc:Fallback></mc:AlternateContent></w:r><w:r w:rsidR="00D6299E" w:rsidRPr="00D6299E"><w:t xml:space="preserve"> </w:t></w:r><w:r w:rsidR="00D6299E"><w:t>N/A</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$rng81.InsertXML($xml81)
Write-Output "Done paragraph 81"
